$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New schedule data for rows 2-13 (Materia, NRC, Profesor, Dia, HoraInicio, HoraFin, Salon)
$data = @(
    @("Redes Inalambricas", 49362, "LOPEZ - MUNOZ MAURO ALBERTO", "L", "13:00", "13:59", "1CCO4/301"),
    @("Redes Inalambricas", 49362, "LOPEZ - MUNOZ MAURO ALBERTO", "M", "13:00", "14:59", "1CCO3/303"),
    @("Redes Inalambricas", 49362, "LOPEZ - MUNOZ MAURO ALBERTO", "V", "13:00", "14:59", "1CCO4/301"),
    @("Tec.de Inteligencia Artificial", 49245, "TECUANHUEHUE - VERA PEDRO", "L", "11:00", "11:59", "1CCO3/114"),
    @("Tec.de Inteligencia Artificial", 49245, "TECUANHUEHUE - VERA PEDRO", "M", "11:00", "12:59", "1CCO3/114"),
    @("Tec.de Inteligencia Artificial", 49245, "TECUANHUEHUE - VERA PEDRO", "V", "11:00", "12:59", "1CCO5/202"),
    @("Teoria de Control", 49190, "HERNANDEZ - AMECA JOSE LUIS", "L", "12:00", "12:59", "1CCO3/310"),
    @("Teoria de Control", 49190, "HERNANDEZ - AMECA JOSE LUIS", "A", "11:00", "12:59", "1CCO3/310"),
    @("Teoria de Control", 49190, "HERNANDEZ - AMECA JOSE LUIS", "J", "11:00", "12:59", "1CCO1/002"),
    @("Vision y Animacion por Comput.", 49971, "JUAREZ - PEREZ SILVESTRE", "L", "10:00", "10:59", "1CCO4/103"),
    @("Vision y Animacion por Comput.", 49971, "JUAREZ - PEREZ SILVESTRE", "A", "09:00", "10:59", "1CCO4/103"),
    @("Vision y Animacion por Comput.", 49971, "JUAREZ - PEREZ SILVESTRE", "J", "09:00", "10:59", "1CCO4/308")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $ws.Cells.Item($row, 7).Value = $entry[6]
    $row++
}

# Remove now-unused rows 14 and 15 (table shrank from 14 data rows to 12)
$ws.Range("A14:G15").Delete()
